$wb = $excel.ActiveWorkbook

# ---- Sheet: Resumen ----
$wsR = $wb.Worksheets.Item("Resumen")
$wsR.Range("B2").Value = "Z3"
$wsR.Range("C2").Value = 555.5323007892746

# ---- Sheet: Solucion ----
$wsS = $wb.Worksheets.Item("Solucion")
$wsS.Range("A2").Value = "Pedido_47"
$wsS.Range("A3").Value = "Pedido_30"
$wsS.Range("A4").Value = "Pedido_5"
$wsS.Range("A5").Value = "Pedido_60"
$wsS.Range("B5").Value = "S011"
$wsS.Range("A6").Value = "Pedido_28"
$wsS.Range("B6").Value = "S051"
$wsS.Range("A7").Value = "Pedido_6"
$wsS.Range("A8").Value = "Pedido_20"
$wsS.Range("A9").Value = "Pedido_2"
$wsS.Range("A10").Value = "Pedido_40"
$wsS.Range("A11").Value = "Pedido_29"
$wsS.Range("A12").Value = "Pedido_49"
$wsS.Range("B12").Value = "S032"
$wsS.Range("A13").Value = "Pedido_3"
$wsS.Range("B13").Value = "S052"
$wsS.Range("A14").Value = "Pedido_7"
$wsS.Range("B14").Value = "S003"
$wsS.Range("A15").Value = "Pedido_16"
$wsS.Range("B15").Value = "S023"
$wsS.Range("A16").Value = "Pedido_24"
$wsS.Range("B16").Value = "S043"
$wsS.Range("A17").Value = "Pedido_38"
$wsS.Range("B17").Value = "S013"
$wsS.Range("A18").Value = "Pedido_8"
$wsS.Range("B18").Value = "S053"
$wsS.Range("A19").Value = "Pedido_44"
$wsS.Range("B19").Value = "S033"
$wsS.Range("A20").Value = "Pedido_51"
$wsS.Range("B20").Value = "S004"
$wsS.Range("A21").Value = "Pedido_54"
$wsS.Range("A22").Value = "Pedido_46"
$wsS.Range("A23").Value = "Pedido_27"
$wsS.Range("B24").Value = "S054"
$wsS.Range("A25").Value = "Pedido_48"
$wsS.Range("B25").Value = "S034"
$wsS.Range("A26").Value = "Pedido_53"
$wsS.Range("B26").Value = "S045"
$wsS.Range("A27").Value = "Pedido_4"
$wsS.Range("B27").Value = "S005"
$wsS.Range("A28").Value = "Pedido_36"
$wsS.Range("B28").Value = "S055"
$wsS.Range("A29").Value = "Pedido_10"
$wsS.Range("B29").Value = "S046"
$wsS.Range("A30").Value = "Pedido_55"
$wsS.Range("B30").Value = "S025"
$wsS.Range("A31").Value = "Pedido_42"
$wsS.Range("B31").Value = "S015"
$wsS.Range("A32").Value = "Pedido_56"
$wsS.Range("B32").Value = "S035"
$wsS.Range("A33").Value = "Pedido_45"
$wsS.Range("B33").Value = "S056"
$wsS.Range("A34").Value = "Pedido_58"
$wsS.Range("B34").Value = "S006"
$wsS.Range("A35").Value = "Pedido_13"
$wsS.Range("B35").Value = "S026"
$wsS.Range("A36").Value = "Pedido_11"
$wsS.Range("B36").Value = "S016"
$wsS.Range("A37").Value = "Pedido_19"
$wsS.Range("B37").Value = "S047"
$wsS.Range("A38").Value = "Pedido_50"
$wsS.Range("B38").Value = "S036"
$wsS.Range("A39").Value = "Pedido_15"
$wsS.Range("B39").Value = "S027"
$wsS.Range("A40").Value = "Pedido_26"
$wsS.Range("B40").Value = "S057"
$wsS.Range("A41").Value = "Pedido_18"
$wsS.Range("B41").Value = "S007"
$wsS.Range("A42").Value = "Pedido_52"
$wsS.Range("B42").Value = "S037"
$wsS.Range("A43").Value = "Pedido_57"
$wsS.Range("B43").Value = "S017"
$wsS.Range("A44").Value = "Pedido_23"
$wsS.Range("B44").Value = "S048"
$wsS.Range("A45").Value = "Pedido_31"
$wsS.Range("B45").Value = "S028"
$wsS.Range("A46").Value = "Pedido_25"
$wsS.Range("B46").Value = "S008"
$wsS.Range("A47").Value = "Pedido_33"
$wsS.Range("A48").Value = "Pedido_9"
$wsS.Range("B48").Value = "S038"
$wsS.Range("A49").Value = "Pedido_39"
$wsS.Range("B49").Value = "S018"
$wsS.Range("A50").Value = "Pedido_17"
$wsS.Range("A51").Value = "Pedido_1"
$wsS.Range("B51").Value = "S029"
$wsS.Range("A52").Value = "Pedido_22"
$wsS.Range("B52").Value = "S059"
$wsS.Range("A53").Value = "Pedido_21"
$wsS.Range("B53").Value = "S009"
$wsS.Range("A54").Value = "Pedido_59"
$wsS.Range("B54").Value = "S039"
$wsS.Range("A55").Value = "Pedido_34"
$wsS.Range("A56").Value = "Pedido_35"
$wsS.Range("B56").Value = "S030"
$wsS.Range("A57").Value = "Pedido_37"
$wsS.Range("B57").Value = "S050"
$wsS.Range("A58").Value = "Pedido_41"
$wsS.Range("B58").Value = "S010"
$wsS.Range("A59").Value = "Pedido_14"
$wsS.Range("B59").Value = "S040"
$wsS.Range("A60").Value = "Pedido_43"
$wsS.Range("B60").Value = "S060"
$wsS.Range("A61").Value = "Pedido_12"
$wsS.Range("B61").Value = "S020"

# ---- Sheet: Metricas ----
$wsM = $wb.Worksheets.Item("Metricas")
$wsM.Range("B2").Value = 551.7115255703319
$wsM.Range("B3").Value = 554.4410530868201
$wsM.Range("B4").Value = 555.5323007892746
